$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("N4").Value = 0
